$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$used = $ws.UsedRange
$used.Replace("globalmean", "g_mean", 1) | Out-Null
$used.Replace("5_imputations_combined", "combination_all", 1) | Out-Null

$ws.Range("A154").Select() | Out-Null
$ws.Range("H163").Select() | Out-Null
